# Mise à jour des textes de la colonne "Matériel / Condition d'accès / Document
# obligatoire" : remplacement de la mention "Présenter un certificat médical"
# par "Questionnaire médical à remplir lors de l'inscription en ligne", et
# repositionnement de la vue de la feuille.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Badminton : Condition d'accès
$ws.Range("D2").Value = "Condition d'accès :`nQuestionnaire médical à remplir lors de l'inscription en ligne de non contre indication à la pratique du badminton (questionnaire médical pour les mineurs)"

# Row 3 - Tir à l'arc : Document obligatoire (note trailing space preserved)
$ws.Range("D3").Value = "Document obligatoire : `nQuestionnaire médical à remplir lors de l'inscription en ligne "

# Row 4 - Qi-Qong : Document obligatoire
$ws.Range("D4").Value = "Document obligatoire : `nQuestionnaire médical à remplir lors de l'inscription en ligne"

# Row 5 - Aviron : Condition d'accès
$ws.Range("D5").Value = "Condition d'accès :`nQuestionnaire médical à remplir lors de l'inscription en ligne de non contre indication à la pratique de l'aviron"

# Row 6 - Canoe-Kayak : Document obligatoire / Savoir nager (was rich text, becomes plain text)
$ws.Range("D6").Value = "`nDocument obligatoire : `nQuestionnaire médical à remplir lors de l'inscription en ligne`nSavoir nager`n"

# Row 7 - Judo : Document obligatoire
$ws.Range("D7").Value = "Document obligatoire : `nQuestionnaire médical à remplir lors de l'inscription en ligne"

# Row 8 - Kick-Boxing : Document obligatoire
$ws.Range("D8").Value = "Document obligatoire : `nQuestionnaire médical à remplir lors de l'inscription en ligne"

# Row 9 - Cross-Training : Document obligatoire
$ws.Range("D9").Value = "Document obligatoire : `nQuestionnaire médical à remplir lors de l'inscription en ligne"

# Row 10 - Sport nature : Document obligatoire
$ws.Range("D10").Value = "Document obligatoire : `nQuestionnaire médical à remplir lors de l'inscription en ligne"

# Row 11 - Paintball : Document obligatoire
$ws.Range("D11").Value = "Document obligatoire : `nQuestionnaire médical à remplir lors de l'inscription en ligne"

# Row 12 - Escalade : reste vide (aucun changement)

# Row 13 - Baby Gym : Document obligatoire
$ws.Range("D13").Value = "Document obligatoire : `nQuestionnaire médical à remplir lors de l'inscription en ligne"

# Row 14 - Gymnastique Artistique : Document obligatoire
$ws.Range("D14").Value = "Document obligatoire : `nQuestionnaire médical à remplir lors de l'inscription en ligne"

# Row 15 - Gymnastique Rythmique : Document obligatoire
$ws.Range("D15").Value = "Document obligatoire : `nQuestionnaire médical à remplir lors de l'inscription en ligne"

# Repositionnement de la vue (topLeftCell A7->A3, sélection B8->D3)
$ws.Range("D3").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
